{"js": "// AB#126 Seps\u00e1n\u00ed funk\u010dn\u00edch a nefunk\u010dn\u00edch po\u017eadavk\u016f pro diskuze s t\u00fdmem\n//\n// 1) \"...je Vajgl Marek, kter\u00fd...\" -> \"...je RNDr. Vajgl Marek, kter\u00fd...\"\n// 2) Drop the \"p\u0159es protokoly SMTP a RabbitMQ (AMQP)\" clause from the\n//    notification micro-service description.\n// 3) Drop \" pomoc\u00ed RabbitMQ\" from the async-processing NFR bullet.\n// 4) Remove the \"Aplikace postavena na mikroslu\u017eb\u00e1ch...\" bullet entirely.\n// 5) Replace the \"Architektura vyu\u017e\u00edv\u00e1 metodiky DDD...\" bullet text with\n//    the new Controller-Service-Repository description.\n// 6) Remove the \"Extern\u00ed ukl\u00e1d\u00e1n\u00ed dokument\u016f...\" bullet entirely.\n// 7) Remove the \"Antivirov\u00e1 kontrola sta\u017een\u00fdch soubor\u016f.\" bullet entirely.\n\nconst body = context.document.body;\n\n// --- 1) Add the \"RNDr.\" title in front of \"Vajgl Marek\" -----------------\nconst vajglResults = body.search(\"Vajgl Marek\", { matchCase: true, matchWholeWord: false });\nvajglResults.load(\"items\");\nawait context.sync();\nif (vajglResults.items.length > 0) {\n  vajglResults.items[0].insertText(\"RNDr. Vajgl Marek\", Word.InsertLocation.replace);\n}\nawait context.sync();\n\n// --- 2) Shorten the e-mail notification sentence -------------------------\nconst notifResults = body.search(\n  \"E-mailov\u00e9 notifikace jsou zas\u00edl\u00e1ny notifika\u010dn\u00ed mikroslu\u017ebou, kter\u00e1 komunikuje asynchronn\u011b p\u0159es protokoly SMTP a RabbitMQ (AMQP).\",\n  { matchCase: true }\n);\nnotifResults.load(\"items\");\nawait context.sync();\nif (notifResults.items.length > 0) {\n  notifResults.items[0].insertText(\n    \"E-mailov\u00e9 notifikace jsou zas\u00edl\u00e1ny notifika\u010dn\u00ed mikroslu\u017ebou, kter\u00e1 komunikuje asynchronn\u011b.\",\n    Word.InsertLocation.replace\n  );\n}\nawait context.sync();\n\n// --- 3) Shorten the async processing NFR bullet ---------------------------\nconst asyncResults = body.search(\n  \"Asynchronn\u00ed zpracov\u00e1n\u00ed komunikace a notifikac\u00ed pomoc\u00ed RabbitMQ.\",\n  { matchCase: true }\n);\nasyncResults.load(\"items\");\nawait context.sync();\nif (asyncResults.items.length > 0) {\n  asyncResults.items[0].insertText(\n    \"Asynchronn\u00ed zpracov\u00e1n\u00ed komunikace a notifikac\u00ed.\",\n    Word.InsertLocation.replace\n  );\n}\nawait context.sync();\n\n// --- 5) Rewrite the architecture bullet under Supportability --------------\nconst archResults = body.search(\n  \"Architektura vyu\u017e\u00edv\u00e1 metodiky DDD, mediator (MediatR), middleware pro zpracov\u00e1n\u00ed chyb.\",\n  { matchCase: true }\n);\narchResults.load(\"items\");\nawait context.sync();\nif (archResults.items.length > 0) {\n  archResults.items[0].insertText(\n    \"Architektura vyu\u017e\u00edv\u00e1 v\u00edcevrstv\u00fd model Controller\u2013Service\u2013Repository a middleware ControllerAdvice pro zpracov\u00e1n\u00ed chyb.\",\n    Word.InsertLocation.replace\n  );\n}\nawait context.sync();\n\n// --- 4), 6), 7) Remove the three obsolete bullets entirely ---------------\nconst paragraphsToRemove = [\n  \"Aplikace postavena na mikroslu\u017eb\u00e1ch, nasazen\u00ed pomoc\u00ed Docker Compose a Podman.\",\n  \"Extern\u00ed ukl\u00e1d\u00e1n\u00ed dokument\u016f (OneDrive, Google Drive), syst\u00e9m pouze poskytuje odkazy.\",\n  \"Antivirov\u00e1 kontrola sta\u017een\u00fdch soubor\u016f.\"\n];\n\nconst allParagraphs = body.paragraphs;\nallParagraphs.load(\"items/text\");\nawait context.sync();\n\nfor (let i = 0; i < allParagraphs.items.length; i++) {\n  const paragraph = allParagraphs.items[i];\n  if (paragraphsToRemove.indexOf(paragraph.text) !== -1) {\n    paragraph.delete();\n  }\n}\nawait context.sync();\n", "ps1": "# AB#126 Seps\u00e1n\u00ed funk\u010dn\u00edch a nefunk\u010dn\u00edch po\u017eadavk\u016f pro diskuze s t\u00fdmem\n#\n# 1) \"...je Vajgl Marek, kter\u00fd...\" -> \"...je RNDr. Vajgl Marek, kter\u00fd...\"\n# 2) Drop the \"p\u0159es protokoly SMTP a RabbitMQ (AMQP)\" clause from the\n#    notification micro-service description.\n# 3) Drop \" pomoc\u00ed RabbitMQ\" from the async-processing NFR bullet.\n# 4) Remove the \"Aplikace postavena na mikroslu\u017eb\u00e1ch...\" bullet entirely.\n# 5) Replace the \"Architektura vyu\u017e\u00edv\u00e1 metodiky DDD...\" bullet text with\n#    the new Controller-Service-Repository description.\n# 6) Remove the \"Extern\u00ed ukl\u00e1d\u00e1n\u00ed dokument\u016f...\" bullet entirely.\n# 7) Remove the \"Antivirov\u00e1 kontrola sta\u017een\u00fdch soubor\u016f.\" bullet entirely.\n\n$d = $word.ActiveDocument\n\n# --- 1) Add the \"RNDr.\" title in front of \"Vajgl Marek\" -------------------\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Execute(\"Vajgl Marek\", $false, $false, $false, $false, $false, $true, 1, $false, \"RNDr. Vajgl Marek\", 2)\n\n# --- 2) Shorten the e-mail notification sentence ---------------------------\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Execute(\n    \"E-mailov\u00e9 notifikace jsou zas\u00edl\u00e1ny notifika\u010dn\u00ed mikroslu\u017ebou, kter\u00e1 komunikuje asynchronn\u011b p\u0159es protokoly SMTP a RabbitMQ (AMQP).\",\n    $false, $false, $false, $false, $false, $true, 1, $false,\n    \"E-mailov\u00e9 notifikace jsou zas\u00edl\u00e1ny notifika\u010dn\u00ed mikroslu\u017ebou, kter\u00e1 komunikuje asynchronn\u011b.\",\n    2\n)\n\n# --- 3) Shorten the async processing NFR bullet -----------------------------\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Execute(\n    \"Asynchronn\u00ed zpracov\u00e1n\u00ed komunikace a notifikac\u00ed pomoc\u00ed RabbitMQ.\",\n    $false, $false, $false, $false, $false, $true, 1, $false,\n    \"Asynchronn\u00ed zpracov\u00e1n\u00ed komunikace a notifikac\u00ed.\",\n    2\n)\n\n# --- 5) Rewrite the architecture bullet under Supportability ---------------\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Execute(\n    \"Architektura vyu\u017e\u00edv\u00e1 metodiky DDD, mediator (MediatR), middleware pro zpracov\u00e1n\u00ed chyb.\",\n    $false, $false, $false, $false, $false, $true, 1, $false,\n    \"Architektura vyu\u017e\u00edv\u00e1 v\u00edcevrstv\u00fd model Controller\u2013Service\u2013Repository a middleware ControllerAdvice pro zpracov\u00e1n\u00ed chyb.\",\n    2\n)\n\n# --- 4), 6), 7) Remove the three obsolete bullets entirely -----------------\n# Walk backwards so deleting a paragraph doesn't disturb the indices of the\n# ones still to be inspected.\n$toRemove = @(\n    \"Aplikace postavena na mikroslu\u017eb\u00e1ch, nasazen\u00ed pomoc\u00ed Docker Compose a Podman.\",\n    \"Extern\u00ed ukl\u00e1d\u00e1n\u00ed dokument\u016f (OneDrive, Google Drive), syst\u00e9m pouze poskytuje odkazy.\",\n    \"Antivirov\u00e1 kontrola sta\u017een\u00fdch soubor\u016f.\"\n)\n\nfor ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {\n    $p = $d.Paragraphs.Item($i)\n    $t = $p.Range.Text.TrimEnd([char]13)\n    if ($toRemove -contains $t) {\n        $p.Range.Delete()\n    }\n}\n"}
